$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.789.03"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.684.82"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D5").Value = "'313.84"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").Value = "'0.3929"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.3966"
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'51.97"
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'1.424"
$ws.Range("E11").Value = "  -5.02%  "
$ws.Range("D12").Value = "'0.08678"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "'25.23"
$ws.Range("E13").Value = "  -4.58%  "
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").Value = "'7.796"
$ws.Range("E15").Value = "  -4.24%  "
$ws.Range("D16").Value = "'0.00001320"
$ws.Range("E16").Value = "  -3.26%  "
$ws.Range("D17").Value = "1.614.35"
$ws.Range("E17").Value = "  -5.30%  "
$ws.Range("E18").Value = "  -3.32%  "
$ws.Range("D19").Value = "'0.07105"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "'20.13"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("D21").Value = "'7.147"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "'14.10"
$ws.Range("E23").Value = "  -2.28%  "
$ws.Range("D24").Value = "24.782.45"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'2.390"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'23.73"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.774"
$ws.Range("E27").Value = "  -8.31%  "
$ws.Range("D28").Value = "'162.14"
$ws.Range("E28").Value = "  -3.36%  "
$ws.Range("D29").Value = "'5.816"
$ws.Range("E29").Value = "  -3.12%  "
$ws.Range("D30").Value = "'149.80"
$ws.Range("E30").Value = "  +2.75%  "
$ws.Range("D31").Value = "'2.606"
$ws.Range("E31").Value = "  +18.72%  "
$ws.Range("D32").Value = "'7.836"
$ws.Range("E32").Value = "  -8.72%  "
$ws.Range("D33").Value = "1.790.82"
$ws.Range("E33").Value = "  -5.34%  "
$ws.Range("D34").Value = "'0.08460"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("D35").Value = "'0.03085"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("D36").Value = "'1.011"
$ws.Range("E36").Value = "  -4.63%  "
$ws.Range("D37").Value = "'6.948"
$ws.Range("E37").Value = "  -4.17%  "
$ws.Range("D38").Value = "'0.2806"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "'0.09537"
$ws.Range("E39").Value = "  +3.17%  "
$ws.Range("D40").Value = "'10.50"
$ws.Range("E40").Value = "  -4.29%  "
$ws.Range("D41").Value = "'0.7953"
$ws.Range("E41").Value = "  -6.34%  "
$ws.Range("D42").Value = "'13.73"
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("D43").Value = "'1.457"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").Value = "'16.71"
$ws.Range("E44").Value = "  -5.66%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.7173"
$ws.Range("E45").Value = "  -4.32%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.584"
$ws.Range("E46").Value = "  -4.85%  "
$ws.Range("D47").Value = "'4.194"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").Value = "'0.08721"
$ws.Range("E48").Value = "  +5.14%  "
$ws.Range("D49").Value = "'1.004"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "'1.336"
$ws.Range("E50").Value = "  -4.40%  "
$ws.Range("D51").Value = "'138.14"
$ws.Range("E51").Value = "  -1.92%  "
